# Applies the Cactuar_Profits market-data refresh (scheduled runner update).
# Writes refreshed price/profit figures into each job sheet (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) of the active workbook, cell by cell, matching the
# canonical OOXML diff for this commit.

$wb = $excel.ActiveWorkbook

# ---- ALC: 48 cell value updates ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 5503.0557
$ws.Range("I38").Value = 1929.3636
$ws.Range("J38").Value = 11118.857
$ws.Range("K38").Value = 5788.0908
$ws.Range("L38").Value = 33356.571
$ws.Range("M38").Value = -5416.0908
$ws.Range("N38").Value = -34100.571
$ws.Range("H51").Value = 7194.2144
$ws.Range("I51").Value = 4925.5
$ws.Range("J51").Value = 8895.75
$ws.Range("K51").Value = 4925.5
$ws.Range("L51").Value = 8895.75
$ws.Range("M51").Value = -4441.5
$ws.Range("N51").Value = -9863.75
$ws.Range("H86").Value = 55556040
$ws.Range("I86").Value = 76923400
$ws.Range("K86").Value = 76923400
$ws.Range("M86").Value = -76922277
$ws.Range("H89").Value = 55556040
$ws.Range("I89").Value = 76923400
$ws.Range("K89").Value = 384617000
$ws.Range("M89").Value = -384611384
$ws.Range("H92").Value = 52631896
$ws.Range("I92").Value = 62500108
$ws.Range("K92").Value = 62500108
$ws.Range("M92").Value = -62498860
$ws.Range("H112").Value = 3574.027
$ws.Range("J112").Value = 3574.027
$ws.Range("L112").Value = 10722.081
$ws.Range("N112").Value = -12938.081
$ws.Range("H116").Value = 37786320
$ws.Range("I116").Value = 141673820
$ws.Range("K116").Value = 141673820
$ws.Range("M116").Value = -141670378
$ws.Range("H132").Value = 126599.16
$ws.Range("I132").Value = 165687.86
$ws.Range("J132").Value = 19105.25
$ws.Range("K132").Value = 497063.58
$ws.Range("L132").Value = 57315.75
$ws.Range("M132").Value = -494533.58
$ws.Range("N132").Value = -62375.75
$ws.Range("H138").Value = 4659.71
$ws.Range("I138").Value = 1281.6285
$ws.Range("J138").Value = 6478.677
$ws.Range("K138").Value = 3844.8855
$ws.Range("L138").Value = 19436.031
$ws.Range("M138").Value = 1295.1145
$ws.Range("N138").Value = -29716.031

# ---- ARM: 42 cell value updates ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17656.19
$ws.Range("I32").Value = 19712.953
$ws.Range("K32").Value = 19712.953
$ws.Range("M32").Value = -19425.953
$ws.Range("H45").Value = 2031.7894
$ws.Range("I45").Value = 1501.3334
$ws.Range("K45").Value = 1501.3334
$ws.Range("M45").Value = -1124.3334
$ws.Range("H61").Value = 7400.3335
$ws.Range("I61").Value = 4300.3335
$ws.Range("J61").Value = 8433.666999999999
$ws.Range("K61").Value = 4300.3335
$ws.Range("L61").Value = 8433.666999999999
$ws.Range("M61").Value = -4088.3335
$ws.Range("N61").Value = -8857.666999999999
$ws.Range("H74").Value = 11820647
$ws.Range("I74").Value = 25001166
$ws.Range("K74").Value = 25001166
$ws.Range("M74").Value = -25000292
$ws.Range("H77").Value = 11820647
$ws.Range("I77").Value = 25001166
$ws.Range("K77").Value = 125005830
$ws.Range("M77").Value = -125001462
$ws.Range("H102").Value = 2499.25
$ws.Range("I102").Value = 2499.25
$ws.Range("K102").Value = 2499.25
$ws.Range("M102").Value = -877.25
$ws.Range("H132").Value = 16000.559
$ws.Range("I132").Value = 23298.885
$ws.Range("K132").Value = 69896.655
$ws.Range("M132").Value = -67366.655
$ws.Range("H133").Value = 81663
$ws.Range("J133").Value = 81663
$ws.Range("L133").Value = 81663
$ws.Range("N133").Value = -86723
$ws.Range("H136").Value = 7400.3335
$ws.Range("I136").Value = 4300.3335
$ws.Range("J136").Value = 8433.666999999999
$ws.Range("K136").Value = 12901.0005
$ws.Range("L136").Value = 25301.001
$ws.Range("M136").Value = -10351.0005
$ws.Range("N136").Value = -30401.001

# ---- BSM: 7 cell value updates ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3329.577
$ws.Range("I134").Value = 2181.5
$ws.Range("J134").Value = 7156.5
$ws.Range("K134").Value = 6544.5
$ws.Range("L134").Value = 21469.5
$ws.Range("M134").Value = -4009.5
$ws.Range("N134").Value = -26539.5

# ---- CRP: 52 cell value updates ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 173.42857
$ws.Range("I7").Value = 92.59999999999999
$ws.Range("J7").Value = 375.5
$ws.Range("K7").Value = 92.59999999999999
$ws.Range("L7").Value = 375.5
$ws.Range("M7").Value = 20.40000000000001
$ws.Range("N7").Value = -601.5
$ws.Range("H31").Value = 21281808
$ws.Range("I31").Value = 45458464
$ws.Range("J31").Value = 6349.8
$ws.Range("K31").Value = 45458464
$ws.Range("L31").Value = 6349.8
$ws.Range("M31").Value = -45458169
$ws.Range("N31").Value = -6939.8
$ws.Range("H34").Value = 21281808
$ws.Range("I34").Value = 45458464
$ws.Range("J34").Value = 6349.8
$ws.Range("K34").Value = 45458464
$ws.Range("L34").Value = 6349.8
$ws.Range("M34").Value = -45458262
$ws.Range("N34").Value = -6753.8
$ws.Range("H51").Value = 48750
$ws.Range("J51").Value = 48750
$ws.Range("L51").Value = 48750
$ws.Range("N51").Value = -50222
$ws.Range("H61").Value = 48750
$ws.Range("J61").Value = 48750
$ws.Range("L61").Value = 48750
$ws.Range("N61").Value = -49446
$ws.Range("H86").Value = 11999.5
$ws.Range("J86").Value = 11999.5
$ws.Range("L86").Value = 11999.5
$ws.Range("N86").Value = -14245.5
$ws.Range("H89").Value = 11999.5
$ws.Range("J89").Value = 11999.5
$ws.Range("L89").Value = 59997.5
$ws.Range("N89").Value = -71229.5
$ws.Range("H99").Value = 20834.533
$ws.Range("J99").Value = 13999.9
$ws.Range("L99").Value = 13999.9
$ws.Range("N99").Value = -16995.9
$ws.Range("H126").Value = 20834.533
$ws.Range("J126").Value = 13999.9
$ws.Range("L126").Value = 41999.7
$ws.Range("N126").Value = -46939.7
$ws.Range("H132").Value = 26673254
$ws.Range("I132").Value = 27780472
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 83341416
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -83338886
$ws.Range("N132").Value = -305060

# ---- CUL: 11 cell value updates ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10378.5625
$ws.Range("I3").Value = 1871.1666
$ws.Range("J3").Value = 15483
$ws.Range("K3").Value = 5613.4998
$ws.Range("L3").Value = 46449
$ws.Range("M3").Value = -5501.4998
$ws.Range("N3").Value = -46673
$ws.Range("H122").Value = 770.3684
$ws.Range("J122").Value = 858.53845
$ws.Range("L122").Value = 7726.84605
$ws.Range("N122").Value = -12626.84605

# ---- GSM: 26 cell value updates ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 761.2857
$ws.Range("I107").Value = 517.2
$ws.Range("J107").Value = 896.8889
$ws.Range("K107").Value = 517.2
$ws.Range("L107").Value = 896.8889
$ws.Range("M107").Value = 1402.8
$ws.Range("N107").Value = -4736.8889
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""
$ws.Range("H122").Value = 233146.86
$ws.Range("I122").Value = 317014.6
$ws.Range("J122").Value = 7349.077
$ws.Range("K122").Value = 951043.7999999999
$ws.Range("L122").Value = 22047.231
$ws.Range("M122").Value = -948593.7999999999
$ws.Range("N122").Value = -26947.231
$ws.Range("H132").Value = 77536.14
$ws.Range("I132").Value = 100970.86
$ws.Range("K132").Value = 302912.58
$ws.Range("M132").Value = -300382.58
$ws.Range("H137").Value = 300044930
$ws.Range("J137").Value = 300044930
$ws.Range("L137").Value = 300044930
$ws.Range("N137").Value = -300055130

# ---- LTW: 39 cell value updates ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1653.7142
$ws.Range("I61").Value = 1636.6471
$ws.Range("K61").Value = 1636.6471
$ws.Range("M61").Value = -1434.6471
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""
$ws.Range("H70").Value = 72482
$ws.Range("J70").Value = 72482
$ws.Range("L70").Value = 72482
$ws.Range("N70").Value = -73022
$ws.Range("H73").Value = 72482
$ws.Range("J73").Value = 72482
$ws.Range("L73").Value = 72482
$ws.Range("N73").Value = -74354
$ws.Range("H93").Value = 405.8889
$ws.Range("I93").Value = 321.85715
$ws.Range("K93").Value = 321.85715
$ws.Range("M93").Value = 926.14285
$ws.Range("H100").Value = 2228.8
$ws.Range("I100").Value = 2254.2222
$ws.Range("K100").Value = 2254.2222
$ws.Range("M100").Value = -1713.2222
$ws.Range("H113").Value = 1653.7142
$ws.Range("I113").Value = 1636.6471
$ws.Range("K113").Value = 1636.6471
$ws.Range("M113").Value = 533.3529000000001
$ws.Range("H132").Value = 2169.7173
$ws.Range("I132").Value = 2200.902
$ws.Range("J132").Value = 1759.8572
$ws.Range("K132").Value = 6602.706
$ws.Range("L132").Value = 5279.571599999999
$ws.Range("M132").Value = -4072.706
$ws.Range("N132").Value = -10339.5716

# ---- WVR: 19 cell value updates ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 782.6
$ws.Range("I113").Value = 720.2778
$ws.Range("J113").Value = 942.8570999999999
$ws.Range("K113").Value = 2160.8334
$ws.Range("L113").Value = 2828.5713
$ws.Range("M113").Value = 9.166600000000017
$ws.Range("N113").Value = -7168.5713
$ws.Range("H126").Value = 4569.1113
$ws.Range("I126").Value = 4396.706
$ws.Range("K126").Value = 13190.118
$ws.Range("M126").Value = -10720.118
$ws.Range("H132").Value = 3240.225
$ws.Range("I132").Value = 1380.6333
$ws.Range("K132").Value = 4141.8999
$ws.Range("M132").Value = -1611.8999
$ws.Range("H136").Value = 7238.6494
$ws.Range("I136").Value = 1328.925
$ws.Range("K136").Value = 3986.775
$ws.Range("M136").Value = -1436.775

Write-Output "Updated 244 cells across 8 sheets"
